$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.680.93'
$ws.Range("E2").Value = '  +6.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.505.10'
$ws.Range("E3").Value = '  +8.36%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.39'
$ws.Range("E5").Value = '  +12.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.82'
$ws.Range("E6").Value = '  +15.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("E8").Value = '  +10.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.501.60'
$ws.Range("E9").Value = '  +8.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0989'
$ws.Range("E10").Value = '  +13.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.54'
$ws.Range("E11").Value = '  +7.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("E12").Value = '  +10.23%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.938.54'
$ws.Range("E14").Value = '  +10.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '55.680.41'
$ws.Range("E15").Value = '  +6.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.87'
$ws.Range("E16").Value = '  +11.28%  '
$ws.Range("E17").Value = '  +17.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.505.64'
$ws.Range("E18").Value = '  +9.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.40'
$ws.Range("E19").Value = '  +12.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.17'
$ws.Range("E20").Value = '  +9.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  +13.42%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +11.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.35'
$ws.Range("E24").Value = '  +10.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +17.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.414'
$ws.Range("E26").Value = '  +14.39%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.623.19'
$ws.Range("E28").Value = '  +9.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  +9.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0796'
$ws.Range("E30").Value = '  +18.95%  '
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.56'
$ws.Range("E32").Value = '  +6.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.30'
$ws.Range("E33").Value = '  +7.91%  '
$ws.Range("E34").Value = '  +13.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.24'
$ws.Range("E35").Value = '  +12.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.879'
$ws.Range("E36").Value = '  +8.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("E37").Value = '  +7.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.13'
$ws.Range("E38").Value = '  +13.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.26'
$ws.Range("E39").Value = '  +7.93%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0560'
$ws.Range("E40").Value = '  +12.64%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.614'
$ws.Range("E41").Value = '  +17.73%  '
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.45'
$ws.Range("E43").Value = '  +10.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.33'
$ws.Range("E44").Value = '  +9.95%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.74'
$ws.Range("E45").Value = '  +22.30%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.003.21'
$ws.Range("E46").Value = '  +6.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0916'
$ws.Range("E47").Value = '  +11.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.14'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '255.69'
$ws.Range("E49").Value = '  +37.21%  '
$ws.Range("E50").Value = '  +10.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.55'
$ws.Range("E51").Value = '  +12.80%  '
